# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "data last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 13:25"

# Row 13 - India
$ws.Range("B13").Value = 158970
$ws.Range("C13").Value = 884
$ws.Range("D13").Value = 67926
$ws.Range("E13").Value = 86503
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 4541

# Row 67 - Finlandia
$ws.Range("D67").Value = 5500
$ws.Range("E67").Value = 930

# Row 77 - Uzbekistan
$ws.Range("D77").Value = 2685
$ws.Range("E77").Value = 697

# Row 90 - Republica de Macedonia
$ws.Range("B90").Value = 2077
$ws.Range("C90").Value = 38
$ws.Range("D90").Value = 1486
$ws.Range("E90").Value = 470
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 121

# Row 102 - Sri Lanka
$ws.Range("B102").Value = 1471
$ws.Range("C102").Value = 2
$ws.Range("E102").Value = 716

# Row 107 - Libano
$ws.Range("B107").Value = 1168
$ws.Range("C107").Value = 7
$ws.Range("E107").Value = 450

# Row 134 - Malta
$ws.Range("B134").Value = 616
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 501
$ws.Range("E134").Value = 108

# Row 162 - Gibraltar
$ws.Range("B162").Value = 158
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 11

$wb.Save()
